# Weekly price update: insert a new observation row for the week of
# 2022-02-08 (serial 44606) ahead of the existing history, pushing the
# prior rows (old 49-71) down by one (new 50-72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 49; this shifts rows 49:71
# down to 50:72 and grows the used range to A1:R72.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Cells.Item(49, 1).Value  = 1
$ws.Cells.Item(49, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(49, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(49, 4).Value  = 44606
$ws.Cells.Item(49, 5).Value  = 15
$ws.Cells.Item(49, 6).Value  = 100112038
$ws.Cells.Item(49, 7).Value  = "Cebollín baby"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Segunda"
$ws.Cells.Item(49, 10).Value = 200
$ws.Cells.Item(49, 11).Value = 2000
$ws.Cells.Item(49, 12).Value = 2500
$ws.Cells.Item(49, 13).Value = 2250
$ws.Cells.Item(49, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(49, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value = 1125
$ws.Cells.Item(49, 17).Value = 2
$ws.Cells.Item(49, 18).Value = "Hortaliza"
